$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5 through 24 all get an explicit custom row height of 14.25 points.
$ws.Range("A5:A24").EntireRow.RowHeight = 14.25
